# Generate Report for Handback
#
# This script updates the localization-status workbook to reflect that the
# de-de and zh-cn handback packages have now been produced:
#  - The "Status" cells for both languages flip from "Ready for handoff" to
#    "Handed back: in sync with en-US".
#  - The "Latest Target File" / "Latest Handback File" columns are populated
#    with the source .md file (as a hyperlink, like column A) and the
#    generated xlf handback file name, for both tracked source files.
#  - The "Latest Handback DateTime" is stamped for each language.
#  - A couple of columns are widened so the new, longer values are readable.

$wb = $excel.ActiveWorkbook

$overviewSheet = $wb.Worksheets.Item("Overview")
$zhSheet       = $wb.Worksheets.Item("zh-cn")
$deSheet       = $wb.Worksheets.Item("de-de")

$statusText = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Status column updates ("Ready for handoff" -> "Handed back: in sync with en-US")
# ---------------------------------------------------------------------
$overviewSheet.Range("E2").Value = $statusText
$overviewSheet.Range("F2").Value = $statusText
$overviewSheet.Range("E3").Value = $statusText
$overviewSheet.Range("F3").Value = $statusText

$zhSheet.Range("C2").Value = $statusText
$zhSheet.Range("C3").Value = $statusText

$deSheet.Range("C2").Value = $statusText
$deSheet.Range("C3").Value = $statusText

# ---------------------------------------------------------------------
# zh-cn sheet: populate Latest Target File / Latest Handback File / DateTime
# ---------------------------------------------------------------------
$zhSheet.Range("I2").Value = "c021d156-42a8-4474-b910-e8fc5b1022ae.md"
$zhSheet.Range("I2").Style = "HyperLink"
$zhSheet.Hyperlinks.Add($zhSheet.Range("I2"), `
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/91c75e635f3a0f90750c137fd56f575349aef38b/e2e/c021d156-42a8-4474-b910-e8fc5b1022ae.md", `
    [Type]::Missing, [Type]::Missing, "c021d156-42a8-4474-b910-e8fc5b1022ae.md") | Out-Null
$zhSheet.Range("J2").Value = "c021d156-42a8-4474-b910-e8fc5b1022ae.1e07495d5293488fb22fe1151aaca6ce53a6155e.zh-cn.xlf"

$zhSheet.Range("I3").Value = "f468ab25-761c-4f62-988a-f4574bb07bb2.md"
$zhSheet.Range("I3").Style = "HyperLink"
$zhSheet.Hyperlinks.Add($zhSheet.Range("I3"), `
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/91c75e635f3a0f90750c137fd56f575349aef38b/e2e/f468ab25-761c-4f62-988a-f4574bb07bb2.md", `
    [Type]::Missing, [Type]::Missing, "f468ab25-761c-4f62-988a-f4574bb07bb2.md") | Out-Null
$zhSheet.Range("J3").Value = "f468ab25-761c-4f62-988a-f4574bb07bb2.d3b1f3d9c295ff252242bf478b69764365483848.zh-cn.xlf"

# Both rows were handed back for zh-cn at the same timestamp.
$zhSheet.Range("K2").Value = "2016-08-16 06:42:05"
$zhSheet.Range("K3").Value = "2016-08-16 06:42:05"

# ---------------------------------------------------------------------
# de-de sheet: populate Latest Target File / Latest Handback File / DateTime
# ---------------------------------------------------------------------
$deSheet.Range("I2").Value = "c021d156-42a8-4474-b910-e8fc5b1022ae.md"
$deSheet.Range("I2").Style = "HyperLink"
$deSheet.Hyperlinks.Add($deSheet.Range("I2"), `
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/91c75e635f3a0f90750c137fd56f575349aef38b/e2e/c021d156-42a8-4474-b910-e8fc5b1022ae.md", `
    [Type]::Missing, [Type]::Missing, "c021d156-42a8-4474-b910-e8fc5b1022ae.md") | Out-Null
$deSheet.Range("J2").Value = "c021d156-42a8-4474-b910-e8fc5b1022ae.1e07495d5293488fb22fe1151aaca6ce53a6155e.de-de.xlf"

$deSheet.Range("I3").Value = "f468ab25-761c-4f62-988a-f4574bb07bb2.md"
$deSheet.Range("I3").Style = "HyperLink"
$deSheet.Hyperlinks.Add($deSheet.Range("I3"), `
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/91c75e635f3a0f90750c137fd56f575349aef38b/e2e/f468ab25-761c-4f62-988a-f4574bb07bb2.md", `
    [Type]::Missing, [Type]::Missing, "f468ab25-761c-4f62-988a-f4574bb07bb2.md") | Out-Null
$deSheet.Range("J3").Value = "f468ab25-761c-4f62-988a-f4574bb07bb2.d3b1f3d9c295ff252242bf478b69764365483848.de-de.xlf"

# Both rows were handed back for de-de at the same (later) timestamp.
$deSheet.Range("K2").Value = "2016-08-16 06:42:14"
$deSheet.Range("K3").Value = "2016-08-16 06:42:14"

# ---------------------------------------------------------------------
# Column widening so the newly-populated columns are readable
# ---------------------------------------------------------------------
# 29.9777047293527 (stored width) corresponds to a character ColumnWidth of
# roughly 29.17 once Excel's fixed +5/6 column-width padding is applied.
$wideColumnWidth = 29.1666666666667

$overviewSheet.Columns.Item(5).ColumnWidth = $wideColumnWidth
$overviewSheet.Columns.Item(6).ColumnWidth = $wideColumnWidth

$zhSheet.Columns.Item(3).ColumnWidth = $wideColumnWidth
$zhSheet.Columns.Item(9).ColumnWidth = 40
$zhSheet.Columns.Item(10).ColumnWidth = 40

$deSheet.Columns.Item(3).ColumnWidth = $wideColumnWidth
$deSheet.Columns.Item(9).ColumnWidth = 40
$deSheet.Columns.Item(10).ColumnWidth = 40
